$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(20.76321629402823, 6.977608749697943, 10.27562151199135, 10.16221892711333, 53.19640514517848, 7.344005520526261, 10.0774362490374, 17.67519768825073),
    @(20.56136534886739, 6.545288562992136, 10.15621449206182, 10.13950400312904, 52.08178299165888, 7.344005520526261, 10.09209268637333, 17.67393970846742),
    @(20.44576969684753, 6.264247296778926, 10.08222054588478, 10.12499934053215, 51.39215742805882, 7.344005520526261, 10.10261598828078, 17.67828721713852),
    @(20.40081852552269, 6.145800879674189, 10.05191071914751, 10.11894723115809, 51.11011740943168, 7.344005520526261, 10.10728735383756, 17.68134529523607),
    @(20.39348614657108, 6.125896313449733, 10.04686872945208, 10.11793371346347, 51.06323286136544, 7.344005520526261, 10.10808615715359, 17.6819307377885),
    @(20.44515467217429, 6.262665753540077, 10.08181239341771, 10.12491829344944, 51.38835742474635, 7.344005520526261, 10.10267743738784, 17.67832325266087),
    @(20.6919228460544, 6.831781386851278, 10.2345988211957, 10.15450187114364, 52.81335284625705, 7.344005520526261, 10.08217331218374, 17.67370094145011),
    @(21.23924707961427, 7.82417764835749, 10.52823790623448, 10.20814250218243, 55.55317911498521, 7.344005520526261, 10.0540684554424, 17.70527014309141),
    @(21.67600274356195, 8.478250637466774, 10.73953937687208, 10.24495939471997, 57.51628383688244, 7.344005520526261, 10.04081115358287, 17.75319598954592),
    @(21.8812916236734, 8.759620414251277, 10.83453680030568, 10.26116114932778, 58.395286151267, 7.344005520526261, 10.03638687428098, 17.78034195012398),
    @(21.959898871699, 8.863852382217969, 10.87033426203834, 10.26721901117615, 58.72587409507355, 7.344005520526261, 10.03494257431575, 17.7913860556454),
    @(21.94293205826422, 8.841507047008177, 10.86263269005848, 10.2659177652627, 58.65478095806834, 7.344005520526261, 10.03524335246427, 17.78897358077524),
    @(21.88774174614233, 8.768242016561469, 10.83748548706565, 10.26166107746948, 58.42253073841937, 7.344005520526261, 10.03626341953428, 17.78123525927684),
    @(21.8540468012027, 8.723063649627303, 10.82205875991609, 10.25904369005465, 58.27996772153682, 7.344005520526261, 10.03691833422571, 17.77659474316795),
    @(21.66271265883942, 8.459537777294475, 10.73330721397121, 10.24388967743694, 57.45853464763982, 7.344005520526261, 10.04113262462804, 17.75152915820444),
    @(21.54697008781179, 8.293737220754009, 10.67856271682712, 10.2344540672909, 56.95083404100089, 7.344005520526261, 10.04412948726086, 17.73751848756593),
    @(21.48102482537227, 8.19685040729993, 10.64697010841527, 10.22897544943084, 56.65751231165109, 7.344005520526261, 10.046004417682, 17.72996324385311),
    @(21.45880700017663, 8.163784419350769, 10.636255795973, 10.22711160167364, 56.55798236487642, 7.344005520526261, 10.0466652053677, 17.72749172045709),
    @(21.55922679485298, 8.311544508446762, 10.68440133562155, 10.2354638260526, 57.00501667418063, 7.344005520526261, 10.0437948154065, 17.7389578812116),
    @(21.90392955207754, 8.78982453069232, 10.84487671910193, 10.26291346073596, 58.49081184494388, 7.344005520526261, 10.03595752947905, 17.78348747928742),
    @(22.13423851302335, 9.088909640996951, 10.94872438651351, 10.28040284695567, 59.44852022281055, 7.344005520526261, 10.03218230362264, 17.81704386016255),
    @(22.01088549948932, 8.930514261058102, 10.89339802329256, 10.27110926235294, 58.93867418625582, 7.344005520526261, 10.03407396541906, 17.79872822334609),
    @(21.55368367291763, 8.303498711652564, 10.68176206647097, 10.23500748227912, 56.9805251649322, 7.344005520526261, 10.04394564712904, 17.73830557467),
    @(21.08479358039362, 7.568917841879745, 10.44952250100252, 10.19409191574797, 54.8195875211778, 7.344005520526261, 10.06037428862314, 17.69238019449607),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rowVals = $data[$i]
    $ws.Range("B$row").Value2 = $rowVals[0]
    $ws.Range("C$row").Value2 = $rowVals[1]
    $ws.Range("D$row").Value2 = $rowVals[2]
    $ws.Range("E$row").Value2 = $rowVals[3]
    $ws.Range("F$row").Value2 = $rowVals[4]
    $ws.Range("H$row").Value2 = $rowVals[5]
    $ws.Range("L$row").Value2 = $rowVals[6]
    $ws.Range("M$row").Value2 = $rowVals[7]
}

Write-Output "Updated loading_percent values for rows 2-25 (B,C,D,E,F,H,L,M)"